$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = 0
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = 1
$ws.Range("F20").Value = -6
$ws.Range("F30").Value = -6
$ws.Range("F34").Value = -1
$ws.Range("F38").Value = -3
$ws.Range("F39").Value = -2
